$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 103
$ws1.Range("F18").Value = 124
$ws1.Range("F23").Value = 11930
$ws1.Range("F24").Value = 11943
$ws1.Range("F27").Value = 245

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 3

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 103
$ws4.Range("F22").Value = 124
$ws4.Range("F27").Value = 11930
$ws4.Range("F28").Value = 11943
$ws4.Range("F31").Value = 245
$ws4.Range("F37").Value = 3
